$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 137, shifting existing rows 137-151 down to 138-152.
$ws.Rows.Item(137).EntireRow.Insert()

# Populate the newly inserted row 137 with the new weekly record.
$ws.Cells.Item(137, 1).Value = 9
$ws.Cells.Item(137, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(137, 3).Value = "Metropolitana"
$ws.Cells.Item(137, 4).Value = 44918
$ws.Cells.Item(137, 5).Value = 13
$ws.Cells.Item(137, 6).Value = 100112022
$ws.Cells.Item(137, 7).Value = "Arveja Verde"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 25
$ws.Cells.Item(137, 11).Value = 26000
$ws.Cells.Item(137, 12).Value = 28000
$ws.Cells.Item(137, 13).Value = 26960
$ws.Cells.Item(137, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(137, 15).Value = "Carahue"
$ws.Cells.Item(137, 16).Value = 1078
$ws.Cells.Item(137, 17).Value = 25
$ws.Cells.Item(137, 18).Value = "Hortaliza"
